$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trends Status")

# Row 2 (Rapid Decline): Current species (no.) 1 -> 2, Current species conclusive (perc.) 16.7 -> 33.3
$ws.Range("C2").Value = 2
$ws.Range("E2").Value = 33.3

# Row 4 (Stable): Current species (no.) 5 -> 4, Current species conclusive (perc.) 83.3 -> 66.7
$ws.Range("C4").Value = 4
$ws.Range("E4").Value = 66.7

$wb.Save()
